# Auto-generated PowerShell Excel COM-interop script
# Adds rows 59-125 to Sheet1 with new English/Japanese vocabulary pairs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('diary', '日記|にっき'),
    @('to fill in', '記入する|きにゅうする'),
    @('(news) article', '記事|きじ'),
    @('to memorize', '暗記する|あんきする'),
    @('bank', '銀行|ぎんこう'),
    @('silver medal', '銀メダル|ぎんメダル'),
    @('land covered with snow', '銀世界|ぎんせかい'),
    @('one time', '一回|いっかい'),
    @('out-of-service bus', '回送バス|かいそうバス'),
    @('last inning; last episode', '最終回|さいしゅうかい'),
    @('to turn', '回す|まわす'),
    @('evening', '夕方|ゆうがた'),
    @('dinner', '夕食|ゆうしょく'),
    @('Tanabata', '七夕|たなばた'),
    @('setting sun', '夕日|ゆうひ'),
    @('evening newspaper', '夕刊|ゆうかん'),
    @('Mr./Ms. Kuroki', '黒木さん|くろきさん'),
    @('black', '黒い|くろい'),
    @('black and white photograph', '白黒写真|しろくろしゃしん'),
    @('blackboard', '黒板|こくばん'),
    @('a thing to take care of', '用事|ようじ'),
    @('to prepare', '用意する|よういする'),
    @('for children', '子供用|こどもよう'),
    @('cost', '費用|ひよう'),
    @('weekend', '週末|しゅうまつ'),
    @('end of the month', '月末|げつまつ'),
    @('year-end', '年末|ねんまつ'),
    @('final examination', '期末試験|きまつしけん'),
    @('the end', '末|すえ'),
    @('to wait', '待つ|まつ'),
    @('waiting room', '待合室|まちあいしつ'),
    @('to expect', '期待する|きたいする'),
    @('invitation', '招待|しょうたい'),
    @('over-time work', '残業|ざんぎょう'),
    @('to leave', '残す|のこす'),
    @('regrettable', '残念|ざんねん'),
    @('regret', '心残り|こころのこり'),
    @('account balance', '残高|ざんだか'),
    @('station', '駅|えき'),
    @('Tokyo Station', '東京駅|とうきょうえき'),
    @('train station attendant', '駅員|えきいん'),
    @('near/in front of the station', '駅前|えきまえ'),
    @('the first', '一番|いちばん'),
    @('number', '番号|ばんごう'),
    @('TV program', '番組|ばんぐみ'),
    @('police box', '交番|こうばん'),
    @('turn; order', '順番|じゅんばん'),
    @('to explain', '説明する|せつめいする'),
    @('novel', '小説|しょうせつ'),
    @('novelist', '小説家|しょうせつか'),
    @('to preach', '説教する|せっきょうする'),
    @('to guide', '案内する|あんないする'),
    @('information desk', '案内所|あんないじょ'),
    @('idea; proposal', '案|あん'),
    @('proposal', '提案|ていあん'),
    @('my wife', '家内|かない'),
    @('domestic', '国内|こくない'),
    @('internal medicine', '内科|ないか'),
    @('inside', '内側|うちがわ'),
    @('to forget', '忘れる|わすれる'),
    @('lost article', '忘れ物|わすれもの'),
    @('year-end party', '忘年会|ぼうねんかい'),
    @('to keep (a promise)', '守る|まもる'),
    @('absence; not at home', '留守|るす'),
    @('answering machine', '留守番電話|るすばんでんわ'),
    @('charm', 'お守り|おまもり'),
    @('security guard', '守衛|しゅえい')
)

$startRow = 59
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
